# Update "想去人数" (column F) values per the commit diff.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 702
$ws1.Range("F3").Value = 56
$ws1.Range("F4").Value = 2009
$ws1.Range("F5").Value = 5882
$ws1.Range("F6").Value = 1661
$ws1.Range("F7").Value = 180
$ws1.Range("F8").Value = 3329
$ws1.Range("F9").Value = 698
$ws1.Range("F10").Value = 50
$ws1.Range("F11").Value = 1395
$ws1.Range("F12").Value = 4638
$ws1.Range("F14").Value = 1747
$ws1.Range("F15").Value = 10
$ws1.Range("F16").Value = 56
$ws1.Range("F17").Value = 61
$ws1.Range("F18").Value = 201
$ws1.Range("F20").Value = 1040
$ws1.Range("F21").Value = 315
$ws1.Range("F22").Value = 84
$ws1.Range("F23").Value = 25
$ws1.Range("F24").Value = 92
$ws1.Range("F26").Value = 218
$ws1.Range("F28").Value = 1134
$ws1.Range("F29").Value = 424
$ws1.Range("F30").Value = 102
$ws1.Range("F31").Value = 219
$ws1.Range("F32").Value = 443
$ws1.Range("F35").Value = 1777
$ws1.Range("F36").Value = 2281
$ws1.Range("F37").Value = 1066
$ws1.Range("F39").Value = 3
$ws1.Range("F40").Value = 286
$ws1.Range("F41").Value = 648
$ws1.Range("F42").Value = 409
$ws1.Range("F43").Value = 49
$ws1.Range("F44").Value = 687
$ws1.Range("F45").Value = 40
$ws1.Range("F46").Value = 455
$ws1.Range("F47").Value = 438

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 798

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 798
$ws4.Range("F3").Value = 702
$ws4.Range("F4").Value = 56
$ws4.Range("F5").Value = 2009
$ws4.Range("F6").Value = 5882
$ws4.Range("F7").Value = 1661
$ws4.Range("F8").Value = 180
$ws4.Range("F9").Value = 3329
$ws4.Range("F10").Value = 50
$ws4.Range("F11").Value = 1395
$ws4.Range("F12").Value = 4638
$ws4.Range("F13").Value = 1747
$ws4.Range("F14").Value = 10
$ws4.Range("F16").Value = 56
$ws4.Range("F19").Value = 61
$ws4.Range("F20").Value = 201
$ws4.Range("F23").Value = 1040
$ws4.Range("F24").Value = 315
$ws4.Range("F25").Value = 92
$ws4.Range("F26").Value = 218
$ws4.Range("F28").Value = 1134
$ws4.Range("F29").Value = 424
$ws4.Range("F30").Value = 102
$ws4.Range("F31").Value = 219
$ws4.Range("F33").Value = 1777
$ws4.Range("F34").Value = 2281
$ws4.Range("F35").Value = 1066
$ws4.Range("F39").Value = 286
$ws4.Range("F40").Value = 648
$ws4.Range("F41").Value = 409
$ws4.Range("F42").Value = 687
$ws4.Range("F43").Value = 455
$ws4.Range("F44").Value = 438
